$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.633.58'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").Value = '3.152.75'
$ws.Range("E3").Value = '  +1.30%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.01'
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.38'
$ws.Range("E6").Value = '  +1.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").Value = '  +11.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.35'

$ws.Range("E10").Value = '  +3.00%  '

$ws.Range("E11").Value = '  +3.34%  '

$ws.Range("D13").Value = '3.693.49'
$ws.Range("E13").Value = '  +1.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.91'
$ws.Range("E14").Value = '  +1.69%  '

$ws.Range("E15").Value = '  +5.60%  '

$ws.Range("D16").Value = '58.666.12'
$ws.Range("E16").Value = '  +1.57%  '

$ws.Range("D17").Value = '3.148.59'
$ws.Range("E17").Value = '  +1.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.20'
$ws.Range("E18").Value = '  +4.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  +3.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.14'
$ws.Range("E20").Value = '  +2.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.51'
$ws.Range("E21").Value = '  +6.25%  '

$ws.Range("E22").Value = '  +1.75%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.59'
$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.514'
$ws.Range("E25").Value = '  +2.10%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.99'
$ws.Range("E28").Value = '  +11.57%  '

$ws.Range("D29").Value = '0.0₃0870'
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.88'
$ws.Range("E30").Value = '  +0.97%  '

$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  +2.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.96'
$ws.Range("E32").Value = '  +3.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.20'
$ws.Range("E33").Value = '  +5.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.17'
$ws.Range("E34").Value = '  +2.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.68'
$ws.Range("E35").Value = '  +0.44%  '

$ws.Range("E36").Value = '  +3.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.37'
$ws.Range("E37").Value = '  +9.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.23'
$ws.Range("E38").Value = '  -2.06%  '

$ws.Range("D39").Value = '2.662.32'
$ws.Range("E39").Value = '  +10.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.68'
$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0685'
$ws.Range("E41").Value = '  +2.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.20'
$ws.Range("E42").Value = '  +4.22%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0287'
$ws.Range("E43").Value = '  +8.94%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.710'
$ws.Range("E44").Value = '  +2.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '38.47'
$ws.Range("E45").Value = '  +4.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").Value = '3.192.21'
$ws.Range("E47").Value = '  +1.25%  '

$ws.Range("E48").Value = '  +13.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.980'
$ws.Range("E49").Value = '  +2.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.19'
$ws.Range("E50").Value = '  +2.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.18'
$ws.Range("E51").Value = '  +2.68%  '
